# Updates crypto price/volume table cells to match the latest scrape.
# Numeric-looking "Price" text (e.g. "218.30") must stay TEXT (same as the
# rest of the column, which uses dotted thousands separators like "26.129.94"
# that are not valid numbers). Excel's default General format would silently
# coerce a plain "218.30" into the float 218.3, losing the trailing zero and
# introducing floating point noise, so for those cells we momentarily force
# Text format, write the value, then restore the original "Normal" style so
# no stray formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "26.129.94"
$ws.Range("D3").Value = "1.655.30"
$ws.Range("E3").Value = "  -0.43%  "
$ws.Range("E4").Value = "  -0.32%  "
Set-TextValue "D5" "218.30"
$ws.Range("E5").Value = "  -0.27%  "
Set-TextValue "D6" "0.5294"
$ws.Range("E6").Value = "  +1.38%  "
Set-TextValue "D7" "1.002"
$ws.Range("E7").Value = "  -0.28%  "
Set-TextValue "D8" "0.2612"
$ws.Range("E8").Value = "  -2.22%  "
Set-TextValue "D9" "0.06335"
$ws.Range("E9").Value = "  +0.06%  "
Set-TextValue "D10" "20.40"
$ws.Range("E10").Value = "  -3.29%  "
Set-TextValue "D11" "0.07760"
$ws.Range("E11").Value = "  +0.46%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D12" "4.498"
$ws.Range("E12").Value = "  +1.51%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.631.02"
$ws.Range("E13").Value = "  -2.00%  "
$ws.Range("E14").Value = "  -0.28%  "
$ws.Range("E15").Value = "  -1.40%  "
Set-TextValue "D16" "65.25"
$ws.Range("E16").Value = "  +0.45%  "
$ws.Range("D17").Value = "26.127.08"
Set-TextValue "D19" "4.540"
$ws.Range("E19").Value = "  -2.58%  "
Set-TextValue "D20" "193.50"
$ws.Range("E20").Value = "  -0.41%  "
Set-TextValue "D21" "10.04"
$ws.Range("E21").Value = "  -1.04%  "
Set-TextValue "D22" "6.001"
$ws.Range("E22").Value = "  -1.36%  "
$ws.Range("E23").Value = "  -0.37%  "
$ws.Range("E24").Value = "  +1.26%  "
Set-TextValue "D25" "0.1241"
$ws.Range("E25").Value = "  +0.00%  "
Set-TextValue "D26" "7.274"
$ws.Range("E26").Value = "  +0.51%  "
$ws.Range("E27").Value = "  +0.11%  "
Set-TextValue "D28" "1.433"
$ws.Range("E28").Value = "  +1.72%  "
Set-TextValue "D29" "0.05940"
$ws.Range("E29").Value = "  -0.59%  "
$ws.Range("E30").Value = "  -0.66%  "
Set-TextValue "D31" "3.512"
$ws.Range("E31").Value = "  -3.26%  "
Set-TextValue "D32" "3.237"
$ws.Range("E32").Value = "  -2.30%  "
Set-TextValue "D33" "1.546"
$ws.Range("E33").Value = "  -5.23%  "
$ws.Range("E34").Value = "  -0.15%  "
Set-TextValue "D35" "0.9457"
$ws.Range("E35").Value = "  -3.48%  "
Set-TextValue "D36" "2.760"
Set-TextValue "D37" "0.5634"
$ws.Range("E37").Value = "  -4.48%  "
Set-TextValue "D38" "0.01609"
$ws.Range("E38").Value = "  +1.01%  "
Set-TextValue "D39" "5.842"
$ws.Range("E39").Value = "  -1.77%  "
Set-TextValue "D40" "0.8468"
$ws.Range("E40").Value = "  -1.61%  "
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue "D42" "101.17"
$ws.Range("E42").Value = "  +1.42%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "1.010.96"
$ws.Range("E43").Value = "  -1.88%  "
$ws.Range("D44").Value = "1.800.42"
$ws.Range("E44").Value = "  -0.21%  "
Set-TextValue "D45" "56.88"
$ws.Range("E45").Value = "  -0.72%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.0₈105"
$ws.Range("E46").Value = "  -5.94%  "
$ws.Range("B47").Value = "Frax"
$ws.Range("C47").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextValue "D47" "1.003"
$ws.Range("E47").Value = "  -0.15%  "
$ws.Range("E48").Value = "  +1.39%  "
$ws.Range("E49").Value = "  -0.70%  "
Set-TextValue "D50" "1.468"
$ws.Range("E50").Value = "  -0.49%  "
Set-TextValue "D51" "7.756"
$ws.Range("E51").Value = "  -3.95%  "
